$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph and the copyright ("... Contact:
# luizeleno@usp.br ...") paragraph that follows it. These two paragraphs, plus
# the blank paragraph immediately after the copyright notice, are removed
# entirely (this matches the corresponding Jekyll build no longer emitting
# that footer block). The blank paragraph preceding "Ver no Jupiter" and the
# page-break paragraph that follows the whole block are left untouched.
$paras = $d.Paragraphs
$n = $paras.Count

$startIndex = -1
$copyIndex = -1

for ($i = 1; $i -le $n; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($startIndex -eq -1 -and $t -like "Ver no Jupiter*") {
        $startIndex = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $copyIndex = $i
    }
}

$startPara = $paras.Item($startIndex)
# also remove the blank paragraph that immediately follows the copyright text
$endPara = $paras.Item($copyIndex + 1)

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
